$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing "notes" column (S), shifting it
# (and nothing else, since it was the last column) to column W.
$ws.Columns("S:V").Insert()

# New header / type rows
$ws.Range("S4").Value = "float"
$ws.Range("T4").Value = "float"
$ws.Range("U4").Value = "string"
$ws.Range("V4").Value = "string"

$ws.Range("S5").Value = "projectileScale"
$ws.Range("T5").Value = "impactScale"
$ws.Range("U5").Value = "fireSfx"
$ws.Range("V5").Value = "impactSfx"

# Runic Revolver
$ws.Range("S6").Value = "0.58"
$ws.Range("T6").Value = "0.9"
$ws.Range("U6").Value = "ui/assets/sfx/weapons/runic_revolver_fire.wav"
$ws.Range("V6").Value = "ui/assets/sfx/weapons/runic_revolver_hit.wav"

# Chorus Ray
$ws.Range("S7").Value = "0.8"
$ws.Range("T7").Value = "1.1"
$ws.Range("U7").Value = "ui/assets/sfx/weapons/chorus_ray_fire.wav"
$ws.Range("V7").Value = "ui/assets/sfx/weapons/chorus_ray_hit.wav"

# Tidebreaker Launcher
$ws.Range("S8").Value = "0.76"
$ws.Range("T8").Value = "1.15"
$ws.Range("U8").Value = "ui/assets/sfx/weapons/tidebreaker_launcher_fire.wav"
$ws.Range("V8").Value = "ui/assets/sfx/weapons/tidebreaker_launcher_hit.wav"

# Pulse Carbine
$ws.Range("S9").Value = "0.6"
$ws.Range("T9").Value = "0.9"
$ws.Range("U9").Value = "ui/assets/sfx/weapons/pulse_carbine_fire.wav"
$ws.Range("V9").Value = "ui/assets/sfx/weapons/pulse_carbine_hit.wav"

# Umbral Scattergun
$ws.Range("S10").Value = "0.72"
$ws.Range("T10").Value = "1.05"
$ws.Range("U10").Value = "ui/assets/sfx/weapons/umbral_scattergun_fire.wav"
$ws.Range("V10").Value = "ui/assets/sfx/weapons/umbral_scattergun_hit.wav"

# Eclipse Javelin
$ws.Range("S11").Value = "0.68"
$ws.Range("T11").Value = "1.18"
$ws.Range("U11").Value = "ui/assets/sfx/weapons/eclipse_javelin_fire.wav"
$ws.Range("V11").Value = "ui/assets/sfx/weapons/eclipse_javelin_hit.wav"
